# Apply the updated cryptocurrency market data to Sheet1.
# Numeric-looking text values in column D are written with a leading
# apostrophe so Excel keeps them as text (matching the workbook's
# original inlineStr / General-format cells) instead of auto-converting
# them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '65.485.47'
$ws.Range('E2').Value = '  -0.31%  '

# Row 3
$ws.Range('D3').Value = '2.654.62'
$ws.Range('E3').Value = '  -0.07%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').Value = '''598.43'
$ws.Range('E5').Value = '  -1.72%  '

# Row 6
$ws.Range('D6').Value = '''156.40'
$ws.Range('E6').Value = '  -0.35%  '

# Row 7
$ws.Range('E7').Value = '  +0.03%  '

# Row 8
$ws.Range('E8').Value = '  +4.31%  '

# Row 9
$ws.Range('D9').Value = '''0.123'
$ws.Range('E9').Value = '  -1.80%  '

# Row 10
$ws.Range('B10').Value = 'Cardano'
$ws.Range('C10').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D10').Value = '''0.400'
$ws.Range('E10').Value = '  -1.17%  '

# Row 11
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Value = '''5.86'
$ws.Range('E11').Value = '  -2.23%  '

# Row 12
$ws.Range('E12').Value = '  -0.42%  '

# Row 13
$ws.Range('D13').Value = '''29.18'
$ws.Range('E13').Value = '  -3.00%  '

# Row 14
$ws.Range('D14').Value = '''0.0000194'
$ws.Range('E14').Value = '  -5.53%  '

# Row 15
$ws.Range('D15').Value = '3.132.93'
$ws.Range('E15').Value = '  -0.10%  '

# Row 16
$ws.Range('D16').Value = '65.312.24'
$ws.Range('E16').Value = '  -0.25%  '

# Row 17
$ws.Range('D17').Value = '2.650.91'
$ws.Range('E17').Value = '  +0.08%  '

# Row 18
$ws.Range('D18').Value = '''12.60'
$ws.Range('E18').Value = '  -0.90%  '

# Row 19
$ws.Range('D19').Value = '''4.78'
$ws.Range('E19').Value = '  -2.84%  '

# Row 20
$ws.Range('E20').Value = '  +1.17%  '

# Row 21
$ws.Range('D21').Value = '''349.60'
$ws.Range('E21').Value = '  -2.95%  '

# Row 22
$ws.Range('E22').Value = '  -0.05%  '

# Row 23
$ws.Range('D23').Value = '''69.25'
$ws.Range('E23').Value = '  -1.41%  '

# Row 24
$ws.Range('D24').Value = '''0.0000109'
$ws.Range('E24').Value = '  +1.88%  '

# Row 25
$ws.Range('D25').Value = '''9.67'
$ws.Range('E25').Value = '  +0.72%  '

# Row 26
$ws.Range('D26').Value = '''1.61'
$ws.Range('E26').Value = '  -5.08%  '

# Row 27
$ws.Range('D27').Value = '''0.166'
$ws.Range('E27').Value = '  -1.38%  '

# Row 28
$ws.Range('D28').Value = '''1.58'
$ws.Range('E28').Value = '  -3.94%  '

# Row 29
$ws.Range('D29').Value = '''8.04'
$ws.Range('E29').Value = '  -1.27%  '

# Row 30
$ws.Range('B30').Value = 'Binance-PegBSC-USD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D30').Value = '''1.00'
$ws.Range('E30').Value = '  +0.06%  '

# Row 31
$ws.Range('B31').Value = 'Bittensor'
$ws.Range('C31').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D31').Value = '''538.35'
$ws.Range('E31').Value = '  +0.36%  '

# Row 32
$ws.Range('D32').Value = '''2.11'
$ws.Range('E32').Value = '  -5.74%  '

# Row 33
$ws.Range('D33').Value = '''1.74'
$ws.Range('E33').Value = '  -3.22%  '

# Row 34
$ws.Range('D34').Value = '''6.49'
$ws.Range('E34').Value = '  +1.03%  '

# Row 35
$ws.Range('D35').Value = '''5.41'
$ws.Range('E35').Value = '  -2.69%  '

# Row 36
$ws.Range('D36').Value = '''0.420'
$ws.Range('E36').Value = '  -2.92%  '

# Row 37
$ws.Range('D37').Value = '''20.29'
$ws.Range('E37').Value = '  -1.96%  '

# Row 38
$ws.Range('B38').Value = 'FirstDigitalUSD'
$ws.Range('C38').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D38').Value = '''1.00'
$ws.Range('E38').Value = '  +0.06%  '

# Row 39
$ws.Range('B39').Value = 'Monero'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D39').Value = '''159.24'
$ws.Range('E39').Value = '  -2.53%  '

# Row 40
$ws.Range('D40').Value = '''1.93'
$ws.Range('E40').Value = '  -3.84%  '

# Row 41
$ws.Range('D41').Value = '''1.00'
$ws.Range('E41').Value = '  +0.02%  '

# Row 42
$ws.Range('D42').Value = '''42.56'
$ws.Range('E42').Value = '  +1.40%  '

# Row 43
$ws.Range('D43').Value = '''164.17'
$ws.Range('E43').Value = '  -0.97%  '

# Row 44
$ws.Range('D44').Value = '''4.06'
$ws.Range('E44').Value = '  -2.52%  '

# Row 45
$ws.Range('D45').Value = '''2.28'
$ws.Range('E45').Value = '  -2.77%  '

# Row 46
$ws.Range('D46').Value = '''0.0605'
$ws.Range('E46').Value = '  -1.33%  '

# Row 47
$ws.Range('D47').Value = '''22.88'
$ws.Range('E47').Value = '  -1.14%  '

# Row 48
$ws.Range('D48').Value = '''0.0258'
$ws.Range('E48').Value = '  -2.75%  '

# Row 49
$ws.Range('D49').Value = '''0.639'
$ws.Range('E49').Value = '  -2.21%  '

# Row 50
$ws.Range('E50').Value = '  +2.66%  '

# Row 51
$ws.Range('D51').Value = '''19.96'
$ws.Range('E51').Value = '  +0.47%  '
